# Apply Valefor Profits market-data refresh (scheduled runner sync)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3444.875
$ws.Range("J76").Value = 4700
$ws.Range("L76").Value = 4700
$ws.Range("N76").Value = -5330
$ws.Range("H79").Value = 3444.875
$ws.Range("J79").Value = 4700
$ws.Range("L79").Value = 4700
$ws.Range("N79").Value = -6884
$ws.Range("H137").Value = 1422
$ws.Range("I137").Value = 911.1905
$ws.Range("J137").Value = 2762.875
$ws.Range("K137").Value = 2733.5715
$ws.Range("L137").Value = 8288.625
$ws.Range("M137").Value = -183.5715
$ws.Range("N137").Value = -13388.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 45
$ws.Range("I4").Value = 45
$ws.Range("K4").Value = 45
$ws.Range("M4").Value = 71
$ws.Range("H9").Value = 30000
$ws.Range("J9").Value = 30000
$ws.Range("L9").Value = 30000
$ws.Range("N9").Value = -30340
$ws.Range("H20").Value = 30000
$ws.Range("J20").Value = 30000
$ws.Range("L20").Value = 30000
$ws.Range("N20").Value = -30540
$ws.Range("H61").Value = 2780.5952
$ws.Range("I61").Value = 2011.8077
$ws.Range("J61").Value = 4029.875
$ws.Range("K61").Value = 2011.8077
$ws.Range("L61").Value = 4029.875
$ws.Range("M61").Value = -1799.8077
$ws.Range("N61").Value = -4453.875
$ws.Range("H74").Value = 12196711
$ws.Range("I74").Value = 17858362
$ws.Range("J74").Value = 2387.8462
$ws.Range("K74").Value = 17858362
$ws.Range("L74").Value = 2387.8462
$ws.Range("M74").Value = -17857488
$ws.Range("N74").Value = -4135.8462
$ws.Range("H77").Value = 12196711
$ws.Range("I77").Value = 17858362
$ws.Range("J77").Value = 2387.8462
$ws.Range("K77").Value = 89291810
$ws.Range("L77").Value = 11939.231
$ws.Range("M77").Value = -89287442
$ws.Range("N77").Value = -20675.231
$ws.Range("H136").Value = 2780.5952
$ws.Range("I136").Value = 2011.8077
$ws.Range("J136").Value = 4029.875
$ws.Range("K136").Value = 6035.4231
$ws.Range("L136").Value = 12089.625
$ws.Range("M136").Value = -3485.4231
$ws.Range("N136").Value = -17189.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13893032
$ws.Range("I31").Value = 26317372
$ws.Range("J31").Value = 7004.9414
$ws.Range("K31").Value = 26317372
$ws.Range("L31").Value = 7004.9414
$ws.Range("M31").Value = -26317077
$ws.Range("N31").Value = -7594.9414
$ws.Range("H34").Value = 13893032
$ws.Range("I34").Value = 26317372
$ws.Range("J34").Value = 7004.9414
$ws.Range("K34").Value = 26317372
$ws.Range("L34").Value = 7004.9414
$ws.Range("M34").Value = -26317170
$ws.Range("N34").Value = -7408.9414
$ws.Range("H42").Value = 29875
$ws.Range("I42").Value = 10000
$ws.Range("J42").Value = 36500
$ws.Range("K42").Value = 10000
$ws.Range("L42").Value = 36500
$ws.Range("M42").Value = -9407
$ws.Range("N42").Value = -37686
$ws.Range("H99").Value = 9322.414000000001
$ws.Range("I99").Value = 9193.714
$ws.Range("J99").Value = 9442.532999999999
$ws.Range("K99").Value = 9193.714
$ws.Range("L99").Value = 9442.532999999999
$ws.Range("M99").Value = -7695.714
$ws.Range("N99").Value = -12438.533
$ws.Range("H126").Value = 9322.414000000001
$ws.Range("I126").Value = 9193.714
$ws.Range("J126").Value = 9442.532999999999
$ws.Range("K126").Value = 27581.142
$ws.Range("L126").Value = 28327.599
$ws.Range("M126").Value = -25111.142
$ws.Range("N126").Value = -33267.599
$ws.Range("H134").Value = 1666.6818
$ws.Range("I134").Value = 1516.6875
$ws.Range("J134").Value = 2066.6667
$ws.Range("K134").Value = 4550.0625
$ws.Range("L134").Value = 6200.000100000001
$ws.Range("M134").Value = -2015.0625
$ws.Range("N134").Value = -11270.0001
$ws.Range("H135").Value = 34653.848
$ws.Range("I135").Value = 40000
$ws.Range("J135").Value = 34208.332
$ws.Range("K135").Value = 40000
$ws.Range("L135").Value = 34208.332
$ws.Range("N135").Value = -44348.332
$ws.Range("M135").Value = -34930

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 553.2857
$ws.Range("I113").Value = 549
$ws.Range("J113").Value = 554
$ws.Range("K113").Value = 1647
$ws.Range("L113").Value = 1662
$ws.Range("M113").Value = 523
$ws.Range("N113").Value = -6002
$ws.Range("H122").Value = 1025.1666
$ws.Range("I122").Value = 843.53845
$ws.Range("J122").Value = 1497.4
$ws.Range("K122").Value = 7591.84605
$ws.Range("L122").Value = 13476.6
$ws.Range("M122").Value = -5141.84605
$ws.Range("N122").Value = -18376.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 153.25
$ws.Range("I2").Value = 133.66667
$ws.Range("J2").Value = 212
$ws.Range("K2").Value = 133.66667
$ws.Range("L2").Value = 212
$ws.Range("M2").Value = -20.66667000000001
$ws.Range("N2").Value = -438
$ws.Range("H62").Value = 26799.6
$ws.Range("J62").Value = 26799.6
$ws.Range("L62").Value = 26799.6
$ws.Range("N62").Value = -28171.6
$ws.Range("H65").Value = 26799.6
$ws.Range("J65").Value = 26799.6
$ws.Range("L65").Value = 80398.79999999999
$ws.Range("N65").Value = -87262.79999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1308.75
$ws.Range("I46").Value = 985
$ws.Range("J46").Value = 1416.6666
$ws.Range("K46").Value = 985
$ws.Range("L46").Value = 1416.6666
$ws.Range("M46").Value = -797
$ws.Range("N46").Value = -1792.6666
$ws.Range("H93").Value = 11953.182
$ws.Range("I93").Value = 18414.166
$ws.Range("J93").Value = 4200
$ws.Range("K93").Value = 18414.166
$ws.Range("L93").Value = 4200
$ws.Range("M93").Value = -17166.166
$ws.Range("N93").Value = -6696

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 7576482
$ws.Range("I100").Value = 12987384
$ws.Range("J100").Value = 1220
$ws.Range("K100").Value = 25974768
$ws.Range("L100").Value = 2440
$ws.Range("M100").Value = -25974227
$ws.Range("N100").Value = -3522
$ws.Range("H126").Value = 1440.8334
$ws.Range("I126").Value = 622.5
$ws.Range("J126").Value = 1850
$ws.Range("K126").Value = 1867.5
$ws.Range("L126").Value = 5550
$ws.Range("M126").Value = 602.5
$ws.Range("N126").Value = -10490
$ws.Range("H136").Value = 2570.138
$ws.Range("I136").Value = 2753.9092
$ws.Range("K136").Value = 8261.7276
$ws.Range("M136").Value = -5711.7276

Write-Output "Updated 165 cells across 7 sheets"